$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45179 -> 45180) for every data row (rows 2 through 116).
$ws.Range("C2:C116").Value = 45180
